$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.231199300764985
$ws.Range("D2").Value = 0.3119802486860817
$ws.Range("E2").Value = 0.277074934566004
$ws.Range("F2").Value = 1.016062852702859
$ws.Range("G2").Value = 0.478399028328127
$ws.Range("H2").Value = 0.6015429175750455
$ws.Range("J2").Value = 0.3242521349348237
$ws.Range("L2").Value = 0.3263770821254184
$ws.Range("N2").Value = 1.544034455689427
$ws.Range("O2").Value = 2.11909832864103
$ws.Range("B3").Value = 1.172498427706131
$ws.Range("D3").Value = 0.3145382607038649
$ws.Range("E3").Value = 0.2787484308133124
$ws.Range("F3").Value = 1.015695634011891
$ws.Range("G3").Value = 0.4706120892634402
$ws.Range("H3").Value = 0.6017140473454106
$ws.Range("J3").Value = 0.3245086020521342
$ws.Range("L3").Value = 0.2915641797759179
$ws.Range("N3").Value = 1.534132498240723
$ws.Range("O3").Value = 2.102651827298217
$ws.Range("B4").Value = 1.136742812036772
$ws.Range("D4").Value = 0.3162042631775228
$ws.Range("E4").Value = 0.2798555459126337
$ws.Range("F4").Value = 1.016077682851176
$ws.Range("G4").Value = 0.4661790110110786
$ws.Range("H4").Value = 0.6021190051582863
$ws.Range("J4").Value = 0.3247499973653092
$ws.Range("L4").Value = 0.2701451504159706
$ws.Range("N4").Value = 1.528629063277918
$ws.Range("O4").Value = 2.093908432855613
$ws.Range("B5").Value = 1.122245591141905
$ws.Range("D5").Value = 0.3169071961043239
$ws.Range("E5").Value = 0.2803267791594735
$ws.Range("F5").Value = 1.016386161353765
$ws.Range("G5").Value = 0.4644599745647611
$ws.Range("H5").Value = 0.6023595011320424
$ws.Range("J5").Value = 0.3248695935038697
$ws.Range("L5").Value = 0.2614064453309766
$ws.Range("N5").Value = 1.526532059359923
$ws.Range("O5").Value = 2.090686105062105
$ws.Range("B6").Value = 1.119842815865752
$ws.Range("D6").Value = 0.3170253697435186
$ws.Range("E6").Value = 0.2804062414171469
$ws.Range("F6").Value = 1.016446612853976
$ws.Range("G6").Value = 0.4641798115634401
$ws.Range("H6").Value = 0.602403995545302
$ws.Range("J6").Value = 0.324890737923841
$ws.Range("L6").Value = 0.2599547907832402
$ws.Range("N6").Value = 1.526192674952455
$ws.Range("O6").Value = 2.090171615879569
$ws.Range("B7").Value = 1.136546998160753
$ws.Range("D7").Value = 0.316213645851616
$ws.Range("E7").Value = 0.279861819762095
$ws.Range("F7").Value = 1.016081224426301
$ws.Range("G7").Value = 0.4661554733690139
$ws.Range("H7").Value = 0.6021219429020164
$ws.Range("J7").Value = 0.3247515241761576
$ws.Range("L7").Value = 0.2700273378262636
$ws.Range("N7").Value = 1.528600191456206
$ws.Range("O7").Value = 2.093863596040819
$ws.Range("B8").Value = 1.210900455982852
$ws.Range("D8").Value = 0.3128424815185729
$ws.Range("E8").Value = 0.277635476950568
$ws.Range("F8").Value = 1.015810155810165
$ws.Range("G8").Value = 0.4756417955499614
$ws.Range("H8").Value = 0.6015397193522318
$ws.Range("J8").Value = 0.3243232193657875
$ws.Range("L8").Value = 0.3143831033444258
$ws.Range("N8").Value = 1.540501067644243
$ws.Range("O8").Value = 2.113146367178018
$ws.Range("B9").Value = 1.358933833196147
$ws.Range("D9").Value = 0.3069864810221432
$ws.Range("E9").Value = 0.2738982631986158
$ws.Range("F9").Value = 1.020099442242298
$ws.Range("G9").Value = 0.4970112425790063
$ws.Range("H9").Value = 0.6027752942084277
$ws.Range("J9").Value = 0.3241439931874837
$ws.Range("L9").Value = 0.4009894459474026
$ws.Range("N9").Value = 1.56838217912896
$ws.Range("O9").Value = 2.161715208345271
$ws.Range("B10").Value = 1.468995005414001
$ws.Range("D10").Value = 0.3031415650416127
$ws.Range("E10").Value = 0.2715319608101119
$ws.Range("F10").Value = 1.026192485369606
$ws.Range("G10").Value = 0.5144067978424971
$ws.Range("H10").Value = 0.6051305155200026
$ws.Range("J10").Value = 0.3244083014396821
$ws.Range("L10").Value = 0.4643589363087699
$ws.Range("N10").Value = 1.591600178760487
$ws.Range("O10").Value = 2.203970627482192
$ws.Range("B11").Value = 1.51933598725185
$ws.Range("D11").Value = 0.3014911791132651
$ws.Range("E11").Value = 0.2705370533891678
$ws.Range("F11").Value = 1.029603859840648
$ws.Range("G11").Value = 0.5226906305193921
$ws.Range("H11").Value = 0.6065159478723956
$ws.Range("J11").Value = 0.3246131090735318
$ws.Range("L11").Value = 0.4931245627733176
$ws.Range("N11").Value = 1.602748993485505
$ws.Range("O11").Value = 2.224624412692123
$ws.Range("B12").Value = 1.538436960937929
$ws.Range("D12").Value = 0.3008803693429991
$ws.Range("E12").Value = 0.2701719684392643
$ws.Range("F12").Value = 1.030987650726175
$ws.Range("G12").Value = 0.5258808864636535
$ws.Range("H12").Value = 0.6070856864354255
$ws.Range("J12").Value = 0.3247027051620819
$ws.Range("L12").Value = 0.5040078546168161
$ws.Range("N12").Value = 1.607054477091168
$ws.Range("O12").Value = 2.23265145945544
$ws.Range("B13").Value = 1.534321559635259
$ws.Range("D13").Value = 0.3010112891440722
$ws.Range("E13").Value = 0.2702500781651196
$ws.Range("F13").Value = 1.03068553587984
$ws.Range("G13").Value = 0.5251914339949053
$ws.Range("H13").Value = 0.6069609777940315
$ws.Range("J13").Value = 0.3246828754942399
$ws.Range("L13").Value = 0.5016643836458741
$ws.Range("N13").Value = 1.606123504253844
$ws.Range("O13").Value = 2.230913532869863
$ws.Range("B14").Value = 1.52090668360097
$ws.Range("D14").Value = 0.3014406439294284
$ws.Range("E14").Value = 0.2705067841901041
$ws.Range("F14").Value = 1.029715861784979
$ws.Range("G14").Value = 0.5229520252237734
$ws.Range("H14").Value = 0.6065619169249317
$ws.Range("J14").Value = 0.324620239594843
$ws.Range("L14").Value = 0.494020135175532
$ws.Range("N14").Value = 1.603101535424514
$ws.Range("O14").Value = 2.225280675736087
$ws.Range("B15").Value = 1.512694583410052
$ws.Range("D15").Value = 0.3017054785595183
$ws.Range("E15").Value = 0.270665541387082
$ws.Range("F15").Value = 1.029133886581221
$ws.Range("G15").Value = 0.5215872722505992
$ws.Range("H15").Value = 0.6063233531334049
$ws.Range("J15").Value = 0.3245834376479735
$ws.Range("L15").Value = 0.4893365352743615
$ws.Range("N15").Value = 1.601261367380815
$ws.Range("O15").Value = 2.221857204511309
$ws.Range("B16").Value = 1.465710489067874
$ws.Range("D16").Value = 0.3032514044281509
$ws.Range("E16").Value = 0.2715986159327493
$ws.Range("F16").Value = 1.025982414910231
$ws.Range("G16").Value = 0.5138728914601103
$ws.Range("H16").Value = 0.6050462889300974
$ws.Range("J16").Value = 0.3243966081919041
$ws.Range("L16").Value = 0.4624777306702299
$ws.Range("N16").Value = 1.590883324621061
$ws.Range("O16").Value = 2.202649665251101
$ws.Range("B17").Value = 1.43695634282227
$ws.Range("D17").Value = 0.3042250304531731
$ws.Range("E17").Value = 0.2721918674815313
$ws.Range("F17").Value = 1.024212912008338
$ws.Range("G17").Value = 0.5092353148047124
$ws.Range("H17").Value = 0.6043432453424202
$ws.Range("J17").Value = 0.3243035687672347
$ws.Range("L17").Value = 0.4459844188715181
$ws.Range("N17").Value = 1.58466652067564
$ws.Range("O17").Value = 2.19123315949102
$ws.Range("B18").Value = 1.420443555095574
$ws.Range("D18").Value = 0.3047943252530887
$ws.Range("E18").Value = 0.2725407681370378
$ws.Range("F18").Value = 1.023255344702633
$ws.Range("G18").Value = 0.506602773225481
$ws.Range("H18").Value = 0.6039684343959948
$ws.Range("J18").Value = 0.3242580225455995
$ws.Range("L18").Value = 0.4364921629859282
$ws.Range("N18").Value = 1.58114605697682
$ws.Range("O18").Value = 2.184801432205234
$ws.Range("B19").Value = 1.414857089312818
$ws.Range("D19").Value = 0.3049886754832674
$ws.Range("E19").Value = 0.2726602204438784
$ws.Range("F19").Value = 1.022941468899646
$ws.Range("G19").Value = 0.5057174266534332
$ws.Range("H19").Value = 0.6038466085780385
$ws.Range("J19").Value = 0.3242439735316367
$ws.Range("L19").Value = 0.433277286544751
$ws.Range("N19").Value = 1.579963604380197
$ws.Range("O19").Value = 2.182646900928034
$ws.Range("B20").Value = 1.440014606980526
$ws.Range("D20").Value = 0.3041204249886302
$ws.Range("E20").Value = 0.2721279206256817
$ws.Range("F20").Value = 1.024395047715828
$ws.Range("G20").Value = 0.5097253833456961
$ws.Range("H20").Value = 0.6044150264395682
$ws.Range("J20").Value = 0.3243126492389905
$ws.Range("L20").Value = 0.4477407597581191
$ws.Range("N20").Value = 1.585322594267637
$ws.Range("O20").Value = 2.192434520074983
$ws.Range("B21").Value = 1.524845939489865
$ws.Range("D21").Value = 0.3013141481232795
$ws.Range("E21").Value = 0.2704310673010557
$ws.Range("F21").Value = 1.029998182330985
$ws.Range("G21").Value = 0.5236083454605733
$ws.Range("H21").Value = 0.6066779070418136
$ws.Range("J21").Value = 0.3246383114838451
$ws.Range("L21").Value = 0.4962657037756344
$ws.Range("N21").Value = 1.603986895824022
$ws.Range("O21").Value = 2.226929593825759
$ws.Range("B22").Value = 1.580508590815384
$ws.Range("D22").Value = 0.2995625771529395
$ws.Range("E22").Value = 0.2693900399698705
$ws.Range("F22").Value = 1.034196268822015
$ws.Range("G22").Value = 0.5329926494252106
$ws.Range("H22").Value = 0.6084197149190231
$ws.Range("J22").Value = 0.324921272216244
$ws.Range("L22").Value = 0.5279230848143186
$ws.Range("N22").Value = 1.616672458928463
$ws.Range("O22").Value = 2.250674317897534
$ws.Range("B23").Value = 1.550780677224054
$ws.Range("D23").Value = 0.3004898870343125
$ws.Range("E23").Value = 0.2699394567824829
$ws.Range("F23").Value = 1.031906616636874
$ws.Range("G23").Value = 0.5279555920233179
$ws.Range("H23").Value = 0.6074660424877294
$ws.Range("J23").Value = 0.3247638746321329
$ws.Range("L23").Value = 0.5110323717299252
$ws.Range("N23").Value = 1.609857575822559
$ws.Range("O23").Value = 2.23789147927306
$ws.Range("B24").Value = 1.438631908472132
$ws.Range("D24").Value = 0.3041676873878032
$ws.Range("E24").Value = 0.2721568065977085
$ws.Range("F24").Value = 1.024312518065884
$ws.Range("G24").Value = 0.5095037185061528
$ws.Range("H24").Value = 0.6043824826897719
$ws.Range("J24").Value = 0.3243085192057933
$ws.Range("L24").Value = 0.4469467492305057
$ws.Range("N24").Value = 1.585025816196733
$ws.Range("O24").Value = 2.191890974490207
$ws.Range("B25").Value = 1.318654353110446
$ws.Range("D25").Value = 0.3084901514158922
$ws.Range("E25").Value = 0.274842373013529
$ws.Range("F25").Value = 1.018422678855202
$ws.Range("G25").Value = 0.4909331296876047
$ws.Range("H25").Value = 0.6021867194194357
$ws.Range("J25").Value = 0.3241224514253247
$ws.Range("L25").Value = 0.3776035616652678
$ws.Range("N25").Value = 1.560356950195029
$ws.Range("O25").Value = 2.147423034456921
